$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Voto: Assegna i voti che vengono caricati sulla piattaforma"
#        -> "Voto: Assegna i voti che vengono caricati su una piattaforma"
#           split across 3 runs: ": ... su" | " una" | " piattaforma"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(": Assegna i voti che vengono caricati sulla piattaforma", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullStart = $rng.Start
$fullEnd = $rng.End

# Replace "sulla" with "su una" (text only, still a single run at this point)
$sub = $d.Range($fullStart, $fullEnd)
$sub.Find.Execute("sulla", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sullaStart = $sub.Start
$sullaEnd = $sub.End
$r = $d.Range($sullaStart, $sullaEnd)
$r.Text = "su una"

# The run text is now ": Assegna i voti che vengono caricati su una piattaforma" (56 chars)
# run1 = chars[0:40] -> ": Assegna i voti che vengono caricati su"
# run2 = chars[40:44] -> " una"
# run3 = chars[44:56] -> " piattaforma"
$split1 = $fullStart + 40
$split2 = $fullStart + 44
$newEnd  = $fullStart + 56

# Force a run break at split2 by toggling a formatting property (same value
# before/after) on the trailing portion, then do the same for the middle
# portion. This produces distinct <w:r> elements without altering the
# visible formatting.
$part3 = $d.Range($split2, $newEnd)
$part3.Font.Size = 12
$part3.Font.Size = 10

$part2 = $d.Range($split1, $split2)
$part2.Font.Size = 12
$part2.Font.Size = 10

# ---------------------------------------------------------------------------
# Change 2: add a bold ": float" run right after "versioneRilascio"
# ---------------------------------------------------------------------------
$rngV = $d.Content
$rngV.Find.Execute("versioneRilascio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$vEnd = $rngV.End

# Pre-split the existing run right at its end (touching only its very last
# character) so the original "versioneRilascio" run is left untouched once
# we append new text to the split-off tail.
$lastCharV = $d.Range($vEnd - 1, $vEnd)
$lastCharV.Font.Size = 12
$lastCharV.Font.Size = 10

$tailV = $d.Range($vEnd - 1, $vEnd)
$tailV.Text = $tailV.Text + ": float"

# Now separate the appended text into its own run.
$newRunV = $d.Range($vEnd, $vEnd + 7)
$newRunV.Font.Size = 12
$newRunV.Font.Size = 10

# ---------------------------------------------------------------------------
# Change 3: add a bold ": String" run right after "orarioAggiornamento"
# ---------------------------------------------------------------------------
$rngO = $d.Content
$rngO.Find.Execute("orarioAggiornamento", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$oEnd = $rngO.End

$lastCharO = $d.Range($oEnd - 1, $oEnd)
$lastCharO.Font.Size = 12
$lastCharO.Font.Size = 10

$tailO = $d.Range($oEnd - 1, $oEnd)
$tailO.Text = $tailO.Text + ": String"

$newRunO = $d.Range($oEnd, $oEnd + 8)
$newRunO.Font.Size = 12
$newRunO.Font.Size = 10
